# Weekly update: insert 3 new daily price rows ("Fruta, Femacal de La Calera - Palta")
# at the top of the data block (old row 813), pushing the existing rows 813-853
# down to 816-856 (dimension grows from A1:T853 to A1:T856).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at 813-815, shifting existing rows 813-853 down to 816-856.
$ws.Rows('813:815').Insert()

# Row 813 - Hass / Primera, semana del 2021-11-09
$ws.Cells.Item(813,1).Value = 3
$ws.Cells.Item(813,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(813,3).Value = 'Coquimbo'
$ws.Cells.Item(813,4).Value = 44509
$ws.Cells.Item(813,5).Value = 5
$ws.Cells.Item(813,6).Value = 'Fruta'
$ws.Cells.Item(813,7).Value = 100106
$ws.Cells.Item(813,8).Value = 'Oleaginosos'
$ws.Cells.Item(813,9).Value = 100106002
$ws.Cells.Item(813,10).Value = 'Palta'
$ws.Cells.Item(813,11).Value = 'Hass'
$ws.Cells.Item(813,12).Value = 'Primera'
$ws.Cells.Item(813,13).Value = 48
$ws.Cells.Item(813,14).Value = 25000
$ws.Cells.Item(813,15).Value = 25000
$ws.Cells.Item(813,16).Value = 25000
$ws.Cells.Item(813,17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(813,18).Value = 'Provincia de Quillota'
$ws.Cells.Item(813,19).Value = 2500
$ws.Cells.Item(813,20).Value = 10

# Row 814 - Hass / Segunda, semana del 2021-11-09
$ws.Cells.Item(814,1).Value = 3
$ws.Cells.Item(814,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(814,3).Value = 'Coquimbo'
$ws.Cells.Item(814,4).Value = 44509
$ws.Cells.Item(814,5).Value = 5
$ws.Cells.Item(814,6).Value = 'Fruta'
$ws.Cells.Item(814,7).Value = 100106
$ws.Cells.Item(814,8).Value = 'Oleaginosos'
$ws.Cells.Item(814,9).Value = 100106002
$ws.Cells.Item(814,10).Value = 'Palta'
$ws.Cells.Item(814,11).Value = 'Hass'
$ws.Cells.Item(814,12).Value = 'Segunda'
$ws.Cells.Item(814,13).Value = 47
$ws.Cells.Item(814,14).Value = 23000
$ws.Cells.Item(814,15).Value = 23000
$ws.Cells.Item(814,16).Value = 23000
$ws.Cells.Item(814,17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(814,18).Value = 'Provincia de Quillota'
$ws.Cells.Item(814,19).Value = 2300
$ws.Cells.Item(814,20).Value = 10

# Row 815 - Hass / Tercera, semana del 2021-11-09
$ws.Cells.Item(815,1).Value = 3
$ws.Cells.Item(815,2).Value = 'Femacal de La Calera'
$ws.Cells.Item(815,3).Value = 'Coquimbo'
$ws.Cells.Item(815,4).Value = 44509
$ws.Cells.Item(815,5).Value = 5
$ws.Cells.Item(815,6).Value = 'Fruta'
$ws.Cells.Item(815,7).Value = 100106
$ws.Cells.Item(815,8).Value = 'Oleaginosos'
$ws.Cells.Item(815,9).Value = 100106002
$ws.Cells.Item(815,10).Value = 'Palta'
$ws.Cells.Item(815,11).Value = 'Hass'
$ws.Cells.Item(815,12).Value = 'Tercera'
$ws.Cells.Item(815,13).Value = 40
$ws.Cells.Item(815,14).Value = 20000
$ws.Cells.Item(815,15).Value = 20000
$ws.Cells.Item(815,16).Value = 20000
$ws.Cells.Item(815,17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(815,18).Value = 'Provincia de Quillota'
$ws.Cells.Item(815,19).Value = 2000
$ws.Cells.Item(815,20).Value = 10
